# Updated cryptos list on Fri Aug 30 11:52:26 UTC 2024 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for the latest cryptos snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.252.00'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '2.508.66'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '538.62'
$ws.Range('E5').Value = '  -1.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.59'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').Value = '2.508.06'
$ws.Range('E9').Value = '  -2.02%  '
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.37'
$ws.Range('E12').Value = '  -4.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.352'
$ws.Range('E13').Value = '  -2.99%  '
$ws.Range('D14').Value = '2.960.99'
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.24'
$ws.Range('E15').Value = '  -2.40%  '
$ws.Range('D16').Value = '59.163.34'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('E17').Value = '  -2.06%  '
$ws.Range('D18').Value = '2.509.73'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.04'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.66'
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.85'
$ws.Range('E23').Value = '  -2.06%  '
$ws.Range('E24').Value = '  +1.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.421'
$ws.Range('E25').Value = '  -4.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.167'
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.77'
$ws.Range('E28').Value = '  -3.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.77'
$ws.Range('E29').Value = '  -2.96%  '
$ws.Range('D30').Value = '0.0₃0773'
$ws.Range('E30').Value = '  -3.82%  '
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '163.11'
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('E34').Value = '  -3.01%  '
$ws.Range('E35').Value = '  -10.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.48'
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.20'
$ws.Range('E37').Value = '  -6.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.59'
$ws.Range('E38').Value = '  -3.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.92'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.66'
$ws.Range('E40').Value = '  -1.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.806'
$ws.Range('E41').Value = '  -5.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.20'
$ws.Range('E42').Value = '  -9.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '279.00'
$ws.Range('E43').Value = '  -9.00%  '
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('E45').Value = '  +0.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.595'
$ws.Range('E46').Value = '  -2.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '124.68'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0936'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('E50').Value = '  -2.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.72'
$ws.Range('E51').Value = '  -4.32%  '
